$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet contains a yearly data table in columns A:R (years 2007-2021 in
# row 3, values in row 4, ratios in row 5). This adds one more year column
# (S) for 2022, continuing the same pattern/formatting as the previous
# column (R).

# Copy the formatting of column R (rows 2-5, the data rows) into the new
# column S so the new cells pick up the same styles (borders, number
# formats, fonts) as the rest of the table.
$ws.Range("R2:R5").Copy() | Out-Null
$ws.Range("S2:S5").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Fill in the new column's data.
$ws.Range("S3").Value = 2022
$ws.Range("S4").Value = 211650
$ws.Range("S5").Value = 2.9794303052841493

# Match the workbook's current selection, which now points at the new
# last-used cell in row 2.
$ws.Range("S2").Select() | Out-Null
